$wb = $excel.ActiveWorkbook

# 1. "Clients" sheet: I3 changes from 1202 to 1201
$wsClients = $wb.Worksheets.Item("Clients")
$wsClients.Range("I3").Value = 1201

# 2. "Episodes" sheet: add T5 = 999 and T6 = 999 (additional_diagnosis column)
$wsEpisodes = $wb.Worksheets.Item("Episodes")
$wsEpisodes.Range("T5").Value = 999
$wsEpisodes.Range("T6").Value = 999

# 3. "K10+" sheet: add F5:O5 = 9 (k10p_item1 .. k10p_item10)
$wsK10 = $wb.Worksheets.Item("K10+")
$wsK10.Range("F5:O5").Value = 9

# 4. "K5" sheet: add F4:J4 = 9 (k5_item1 .. k5_item5)
$wsK5 = $wb.Worksheets.Item("K5")
$wsK5.Range("F4:J4").Value = 9

# 5. "SDQ" sheet: add G4:AV4 = 9 (sdq_item1 .. sdq_item42)
$wsSdq = $wb.Worksheets.Item("SDQ")
$wsSdq.Range("G4:AV4").Value = 9
